$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns I ("I0") and J ("IF"),
# reusing the same bold/centered/bordered style as the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the I and J data columns for rows 2-75.
$ijData = @(
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 10),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(7, 8),
    @(7, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(6, 6),
    @(7, 8),
    @(7, 7),
    @(6, 7),
    @(6, 7),
    @(5, 6),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(7, 7),
    @(7, 7),
    @(10, 10),
    @(5, 6),
    @(5, 5),
    @(7, 7),
    @(6, 7),
    @(7, 8),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(6, 7),
    @(5, 5),
    @(6, 7),
    @(6, 6),
    @(7, 7),
    @(7, 7),
    @(7, 8),
    @(4, 5),
    @(7, 8),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(5, 6),
    @(5, 5),
    @(5, 6),
    @(7, 7)
)

for ($idx = 0; $idx -lt $ijData.Count; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $ijData[$idx][0]
    $ws.Cells.Item($row, 10).Value = $ijData[$idx][1]
}
